# started working on cuNN. need to read Topology & Training sections
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear cells whose content no longer applies at their current position ---
# Old row 16 (cuNN) content is being replaced / relocated
$ws.Range("C16:G16").ClearContents()
# Old row 17 (Engine) shrinks down to just the MyAlgebra row's 3 columns
$ws.Range("F17:G17").ClearContents()
# Old row 19 (Forecaster) content moves down to row 21
$ws.Range("C19:G19").ClearContents()

# --- New row 14: Core's dependencies ---
$ws.Range("C14").Value = "Core"
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = "Utils"
$ws.Range("F14").Value = "Debugger"
$ws.Range("G14").Value = "ParamMgr"

# --- Row 15: Engine (now with more dependencies, work in progress on cuNN) ---
$ws.Range("C15").Value = "Engine"
$ws.Range("D15").Value = 6
$ws.Range("E15").Value = "Debugger"
$ws.Range("F15").Value = "ParamMgr"
$ws.Range("G15").Value = "Data"
$ws.Range("H15").Value = "Core"
$ws.Range("I15").Value = "cuNN"
$ws.Range("J15").Value = "cuSVM"
$ws.Range("K15").Value = [char]0x2026

# --- Row 17: MyAlgebra (moved down from row 15) ---
$ws.Range("C17").Value = "MyAlgebra"
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = "MyCU"

# --- Row 18: cuNN (moved down from row 16, gained a new dependency) ---
$ws.Range("C18").Value = "cuNN"
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = "TimeSerie"
$ws.Range("F18").Value = "Core"
$ws.Range("G18").Value = "MyAlgebra"
$ws.Range("H18").Value = "Data"

# --- Row 20: Logger (moved down from row 18) ---
$ws.Range("C20").Value = "Logger"
$ws.Range("E20").Value = "DataSource"
$ws.Range("F20").Value = "OraUtils"

# --- Row 21: Forecaster (moved down from row 19) ---
$ws.Range("C21").Value = "Forecaster"
$ws.Range("D21").Value = 3
$ws.Range("E21").Value = "Data"
$ws.Range("F21").Value = "Engine"
$ws.Range("G21").Value = "Logger"

# --- Update selection to match the author's final cursor position ---
$ws.Range("I18").Select()
